$wb = $excel.ActiveWorkbook

# --- Output_flows sheet updates ---
$ws1 = $wb.Worksheets.Item("Output_flows")
$ws1.Range("C2").Value = [double]"1.256328721605063E-14"
$ws1.Range("E2").Value = [double]"6.433968532939948E-12"
$ws1.Range("G2").Value = [double]"3.020832027309992E-12"
$ws1.Range("I2").Value = [double]"3.345957855737681E-14"
$ws1.Range("K2").Value = [double]"5.047917983504061E-16"
$ws1.Range("C7").Value = [double]"1.824986969640542E-06"
$ws1.Range("F7").Value = [double]"9.46412205480233E-12"
$ws1.Range("G7").Value = [double]"0.004388166084645646"
$ws1.Range("I7").Value = [double]"4.860455215802391E-05"
$ws1.Range("J7").Value = [double]"1.168541664970464E-05"
$ws1.Range("C12").Value = [double]"1.042267987870335E-12"
$ws1.Range("E12").Value = [double]"1.183509472265775E-06"
$ws1.Range("I12").Value = [double]"1.387926862542485E-11"
$ws1.Range("J12").Value = [double]"7.230539365951427E-10"
$ws1.Range("C13").Value = [double]"1.146548775147916E-13"
$ws1.Range("D13").Value = [double]"8.49719297282473E-14"
$ws1.Range("E13").Value = [double]"4.855558253396559E-07"
$ws1.Range("I13").Value = [double]"1.52679144208826E-12"
$ws1.Range("J13").Value = [double]"1.555650373312833E-10"
$ws1.Range("C14").Value = [double]"1.340545739457292E-14"
$ws1.Range("D14").Value = [double]"3.973969911780541E-14"
$ws1.Range("E14").Value = [double]"4.869263946896141E-06"
$ws1.Range("I14").Value = [double]"1.785125767952804E-13"
$ws1.Range("J14").Value = [double]"1.055062157244061E-10"
$ws1.Range("C15").Value = [double]"1.289102725666846E-37"
$ws1.Range("D15").Value = [double]"3.248249574839597E-36"
$ws1.Range("E15").Value = [double]"4.270577847930307E-26"
$ws1.Range("I15").Value = [double]"1.716622137830005E-36"
$ws1.Range("K15").Value = [double]"1.3219918090551E-32"
$ws1.Range("C17").Value = [double]"0.0002547208437736293"
$ws1.Range("F17").Value = [double]"1.157299454546695E-06"
$ws1.Range("I17").Value = [double]"0.0135678695167558"
$ws1.Range("J17").Value = [double]"0.7070113045943595"
$ws1.Range("C18").Value = [double]"2.844516828399762E-05"
$ws1.Range("D18").Value = [double]"8.43250887382852E-07"
$ws1.Range("F18").Value = [double]"4.81906294620537E-07"
$ws1.Range("I18").Value = [double]"0.001515150177511291"
$ws1.Range("J18").Value = [double]"0.1543930611725542"
$ws1.Range("C19").Value = [double]"3.351199054811178E-06"
$ws1.Range("D19").Value = [double]"3.973774676743994E-07"
$ws1.Range("F19").Value = [double]"4.869045578697409E-06"
$ws1.Range("I19").Value = [double]"0.0001785037723130402"
$ws1.Range("J19").Value = [double]"0.1055014646903124"
$ws1.Range("C20").Value = [double]"3.222756939643565E-29"
$ws1.Range("D20").Value = [double]"3.248249702328896E-29"
$ws1.Range("F20").Value = [double]"4.270577848261607E-26"
$ws1.Range("I20").Value = [double]"1.716622204665843E-27"
$ws1.Range("K20").Value = [double]"1.321991806593662E-23"

# --- Input_flows sheet updates ---
$ws2 = $wb.Worksheets.Item("Input_flows")
$ws2.Range("C7").Value = [double]"0.004410849679081186"
$ws2.Range("C12").Value = [double]"2.692645902759782E-08"
$ws2.Range("C13").Value = [double]"3.804524440595465E-09"
$ws2.Range("C14").Value = [double]"3.237892007881076E-10"
$ws2.Range("C17").Value = [double]"0.698051099982766"
$ws2.Range("C18").Value = [double]"0.1536599911639926"
$ws2.Range("C19").Value = [double]"0.105366846483339"
$ws2.Range("E22").Value = [double]"6.228524882861391E-19"
$ws2.Range("F42").Value = [double]"3.720554016690005E-14"
$ws2.Range("F47").Value = [double]"3.943136437198241E-05"
$ws2.Range("F52").Value = [double]"1.840849088433856E-11"
$ws2.Range("F53").Value = [double]"2.277864254506563E-12"
$ws2.Range("F54").Value = [double]"3.168714020126818E-13"
$ws2.Range("F55").Value = [double]"1.322169888430996E-32"
$ws2.Range("F57").Value = [double]"0.01839356073783827"
$ws2.Range("F58").Value = [double]"0.002277306266979654"
$ws2.Range("F59").Value = [double]"0.0003168703374407138"
$ws2.Range("F60").Value = [double]"1.321991806593662E-23"
